$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (2015): B7/C7/D7 were the shared string "x" -> now numeric 1
$ws.Range("B7").Value2 = 1
$ws.Range("C7").Value2 = 1
$ws.Range("D7").Value2 = 1

# Row 8 (2016): left block flips to 0, right block (G/H/I) gains 1s
$ws.Range("B8").Value2 = 0
$ws.Range("C8").Value2 = 0
$ws.Range("D8").Value2 = 0
$ws.Range("G8").Value2 = 1
$ws.Range("H8").Value2 = 1
$ws.Range("I8").Value2 = 1

# Row 9 (2017): B9/D9 flip to 0, C9 stays 1; right block gains G9/I9 (no H9)
$ws.Range("B9").Value2 = 0
$ws.Range("C9").Value2 = 1
$ws.Range("D9").Value2 = 0
$ws.Range("G9").Value2 = 1
$ws.Range("I9").Value2 = 1

# Row 10 (2018): D10 flips to 0, B10/C10 unchanged; right block gains I10 only
$ws.Range("B10").Value2 = 1
$ws.Range("C10").Value2 = 1
$ws.Range("D10").Value2 = 0
$ws.Range("I10").Value2 = 1

# Move the selection like the author's saved session
$ws.Range("K9").Select()
